$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the financial/sprint indicators at end of sprint
$ws.Range("B1").Value = 163
$ws.Range("B2").Value = 0.47
$ws.Range("B3").Value = 89
